# Apply cryptos list price/volume update (commit: Updated cryptos list on Wed Nov 20 10:34:23 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '93.393.90'
$ws.Cells.Item(2, 5).Value = '  +1.83%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.118.17'
$ws.Cells.Item(3, 5).Value = '  +0.03%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''237.97'
$ws.Cells.Item(5, 5).Value = '  -3.08%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''615.12'
$ws.Cells.Item(6, 5).Value = '  -0.37%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.32%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +1.69%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''1.00'
$ws.Cells.Item(9, 5).Value = '  -0.04%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.838'
$ws.Cells.Item(10, 5).Value = '  +13.30%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '3.115.89'
$ws.Cells.Item(11, 5).Value = '  +0.03%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -2.46%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'ShibaInu'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(13, 4).Value = '''0.0000245'
$ws.Cells.Item(13, 5).Value = '  -2.02%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).Value = '''35.29'
$ws.Cells.Item(14, 5).Value = '  +1.28%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '93.144.25'
$ws.Cells.Item(15, 5).Value = '  +1.71%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -2.91%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.699.92'
$ws.Cells.Item(17, 5).Value = '  +0.07%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.121.27'
$ws.Cells.Item(18, 5).Value = '  +0.53%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +1.22%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''14.88'
$ws.Cells.Item(20, 5).Value = '  +0.78%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''6.03'
$ws.Cells.Item(21, 5).Value = '  +4.07%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''0.0000201'
$ws.Cells.Item(22, 5).Value = '  -0.81%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''443.42'
$ws.Cells.Item(23, 5).Value = '  -0.64%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''9.11'
$ws.Cells.Item(24, 5).Value = '  -3.15%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'LEO'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(25, 4).Value = '''8.23'
$ws.Cells.Item(25, 5).Value = '  +5.07%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''5.71'
$ws.Cells.Item(26, 5).Value = '  -2.27%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Aptos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(27, 4).Value = '''13.01'
$ws.Cells.Item(27, 5).Value = '  +11.01%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Litecoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(28, 4).Value = '''85.91'
$ws.Cells.Item(28, 5).Value = '  -2.62%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.05%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''0.182'
$ws.Cells.Item(30, 5).Value = '  +9.25%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +1.83%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''0.125'
$ws.Cells.Item(32, 5).Value = '  -11.30%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''9.31'
$ws.Cells.Item(33, 5).Value = '  -0.62%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''1.02'
$ws.Cells.Item(34, 5).Value = '  +2.49%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''8.00'
$ws.Cells.Item(35, 5).Value = '  +3.55%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''0.162'
$ws.Cells.Item(36, 5).Value = '  -8.42%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''25.97'
$ws.Cells.Item(37, 5).Value = '  -0.86%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'PancakeSwap'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(38, 4).Value = '''1.91'
$ws.Cells.Item(38, 5).Value = '  -1.63%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'MantraDAO'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(39, 4).Value = '''3.89'
$ws.Cells.Item(39, 5).Value = '  -7.46%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''0.449'
$ws.Cells.Item(40, 5).Value = '  +2.51%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.12%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''480.54'
$ws.Cells.Item(42, 5).Value = '  -1.98%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''24.00'
$ws.Cells.Item(43, 5).Value = '  +8.12%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''3.32'
$ws.Cells.Item(44, 5).Value = '  -3.13%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.04%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''159.52'
$ws.Cells.Item(46, 5).Value = '  -0.33%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.703'
$ws.Cells.Item(47, 5).Value = '  +0.89%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''1.88'
$ws.Cells.Item(48, 5).Value = '  -1.23%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''1.34'
$ws.Cells.Item(49, 5).Value = '  -0.27%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''4.42'
$ws.Cells.Item(50, 5).Value = '  +0.90%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.27%  '
